$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> D (price) and E (volume) new values.
# Only rows/columns actually changed by the diff are included.
$updates = @{
    2  = @{ D = "42.234.22";  E = "  -0.71%  " }
    3  = @{ D = "2.235.83";   E = "  -0.22%  " }
    4  = @{             E = "  +0.07%  " }
    5  = @{ D = "243.10";     E = "  -0.83%  " }
    6  = @{ D = "0.627";      E = "  -0.21%  " }
    7  = @{ D = "74.94";      E = "  -0.20%  " }
    8  = @{             E = "  +0.08%  " }
    9  = @{             E = "  -2.63%  " }
    10 = @{ D = "42.51";      E = "  -2.80%  " }
    11 = @{             E = "  +1.18%  " }
    12 = @{ D = "6.98";       E = "  -2.88%  " }
    13 = @{             E = "  +0.32%  " }
    14 = @{ D = "2.574.04";   E = "  +0.12%  " }
    15 = @{             E = "  -1.00%  " }
    16 = @{             E = "  -2.03%  " }
    17 = @{ D = "2.248.17";   E = "  +0.68%  " }
    18 = @{ D = "42.118.70";  E = "  -0.38%  " }
    19 = @{ D = "0.0000105";  E = "  +0.97%  " }
    20 = @{ D = "6.24";       E = "  +0.84%  " }
    21 = @{ D = "72.89";      E = "  +1.25%  " }
    22 = @{ D = "11.50";      E = "  +8.18%  " }
    23 = @{ D = "230.54";     E = "  -0.30%  " }
    24 = @{             E = "  -6.01%  " }
    25 = @{ D = "0.999";      E = "  -0.21%  " }
    26 = @{             E = "  -2.85%  " }
    28 = @{             E = "  -0.83%  " }
    29 = @{             E = "  -2.57%  " }
    30 = @{ D = "167.52";     E = "  +0.40%  " }
    31 = @{             E = "  -0.27%  " }
    32 = @{ D = "5.66";       E = "  -5.39%  " }
    33 = @{             E = "  -1.48%  " }
    34 = @{ D = "29.96";      E = "  -0.23%  " }
    35 = @{             E = "  -0.41%  " }
    36 = @{             E = "  -7.37%  " }
    37 = @{ D = "4.36";       E = "  -5.36%  " }
    38 = @{ D = "0.0305";     E = "  -2.38%  " }
    39 = @{ D = "13.32";      E = "  -0.80%  " }
    40 = @{             E = "  -1.47%  " }
    41 = @{ D = "5.72";       E = "  +0.13%  " }
    42 = @{ D = "65.14";      E = "  +2.45%  " }
    43 = @{             E = "  -0.31%  " }
    44 = @{ D = "8.75" }
    45 = @{ D = "104.80";     E = "  -1.65%  " }
    46 = @{             E = "  -2.34%  " }
    47 = @{             E = "  -0.92%  " }
    48 = @{             E = "  -2.21%  " }
    49 = @{             E = "  -0.62%  " }
    50 = @{             E = "  -1.45%  " }
    51 = @{ D = "2.441.99";   E = "  -0.30%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["E"]
    }
}
